# issue #5: stock data output to json file
#
# Adds a "property_category" column to the 股票 (stock) worksheet, filled
# with the literal value "stock" for every data row, and normalizes the
# malformed "6，250" (full-width comma) total value in the last row to a
# plain "6250" string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (H) - this shifts
# the old H:J (date / legislator_name / legislator_id) columns one to the
# right, to I:K, and leaves a blank column H ready for the new field.
$ws.Range("H1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("H1").Value = "property_category"

# Every stock-holding row (rows 2-6) belongs to the "stock" category.
$ws.Range("H2:H6").Value = "stock"

# Fix the malformed total amount for the last row ("台紙"): the scraped
# value used a full-width comma as a thousands separator ("6，250"); store
# the corrected plain digits as text so it matches the other total values'
# semantics used downstream (it stays a text/string cell, not a number).
$ws.Range("G6").Value = "'6250"
